$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(32, 2.5, 0, 0.5100363294283549, 0.1379436731338501),
    @(33, 2.5, 0, 0.1122070988019307, 0.05105317036310832),
    @(34, 2.5, 0, 0.1207146445910136, 0.1465730428695679),
    @(35, 2.5, 0, 0.04501304229100545, 0.03799559672673544),
    @(36, 2.5, 1, 0.05065822601318359, 0.08503950834274292),
    @(37, 2.5, 9, 0.3982537905375163, 1.038130478064219),
    @(38, 2.5, 10, 0.2209033966064453, 0.5444023569424947),
    @(39, 2.5, 10, 0.1254425287246704, 0.5097094178199768),
    @(40, 2.5, 0, 0.05235658486684163, 0.03746788899103801),
    @(41, 2.5, 10, 0.1752272129058838, 0.4498495141665141),
    @(42, 2.5, 8, 0.2288456320762634, 0.4461415926615397),
    @(43, 2.5, 0, 0.05925114552179973, 0.05913585424423218),
    @(44, 2.5, 0, 0.02632412513097127, 0.06713658968607585),
    @(45, 2.5, 0, 0.02119371891021728, 0.03664478858311971),
    @(46, 2.5, 0, 0.0660416046778361, 0.03993319670359294),
    @(47, 2.5, 0, 0.01863482395807902, 0.02957114775975545),
    @(48, 2.5, 0, 0.01867716312408447, 0.04112821420033773),
    @(49, 2.5, 0, 0.05096803903579712, 0.05055300792058309),
    @(50, 2.5, 6, 0.5126078685124715, 0.425563375155131),
    @(51, 2.5, 5, 2.10612739721934, 0.2804775913556417),
    @(52, 2.5, 3, [double]"1.800457636515299e-05", 2.263049674034119)
)

$startRow = 65
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
